$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate new shared strings in the order they must appear in the
# shared string table: accounts(29), vector_col(30), dict_col(31), vector_feature(32)
$ws.Range("B19").Value = "accounts"
$ws.Range("A20").Value = "vector_col"
$ws.Range("A19").Value = "dict_col"
$ws.Range("B20").Value = "vector_feature"

$ws.Range("F13").Select()
